# Insert a new data row above current row 44 (this shifts rows 44..69 down
# to 45..70, growing the sheet's used range from A1:R69 to A1:R70), then
# populate the newly inserted row 44 with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44460
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = 100112022
$ws.Range("G44").Value = "Arveja Verde"
$ws.Range("H44").Value = "Perfection"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 40
$ws.Range("K44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("M44").Value = 35000
$ws.Range("N44").Value = "`$/malla 25 kilos"
$ws.Range("O44").Value = "Provincia de Huasco"
$ws.Range("P44").Value = 1400
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
